$d = $word.ActiveDocument

# Locate the "Joaozinho" paragraph and the "Cleber" paragraph that follows
# it, then delete the full range spanning both paragraphs (text + their
# paragraph marks). This removes the two <w:p> elements entirely while
# leaving the neighbouring paragraphs (" Otávio " before, "Ideias de
# projeto" after) untouched.
$paras = $d.Paragraphs

$start = $null
$end = $null

for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $text = $p.Range.Text

    if ($text -eq "Joaozinho`r") {
        $start = $p.Range.Start
        $next = $paras.Item($i + 1)
        if ($next.Range.Text -eq "Cleber`r") {
            $end = $next.Range.End
        } else {
            $end = $p.Range.End
        }
        break
    }
}

if ($start -ne $null -and $end -ne $null) {
    $delRange = $d.Range($start, $end)
    $delRange.Delete()
}
